$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from A18 (an existing date cell) onto A19/A20
$ws.Range("A18").Copy()
$ws.Range("A19:A20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A19").Value = 45857
$ws.Range("B19").Value = "Besprechung mit Michael bei ihm Zuhause. Basis Fargen stellen und Status zeigen"
$ws.Range("G19").Value = 1

$ws.Range("A20").Value = 45857
$ws.Range("B20").Value = "Besprechungsprotokoll vom 19.07. geschrieben und Email an allen gesendet"
$ws.Range("G20").Value = 1.5

$win = $excel.ActiveWindow
$win.Zoom = 160
$ws.Range("G21").Select()
